$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '34.563.75'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.48%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.801.53'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.04%  '
$ws.Range('E4').Value = '  +0.20%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '224.59'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.41%  '
$ws.Range('E6').Value = '  +0.22%  '
$ws.Range('E7').Value = '  +0.22%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '41.56'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +14.92%  '
$ws.Range('E9').Value = '  -0.15%  '
$ws.Range('E10').Value = '  -1.56%  '
$ws.Range('E11').Value = '  +3.48%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.061.78'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.02%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.799.27'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.09%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '10.92'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.68%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '34.495.04'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.33%  '
$ws.Range('E16').Value = '  -0.51%  '
$ws.Range('E17').Value = '  -0.54%  '
$ws.Range('E18').Value = '  -2.18%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '240.69'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.27%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0767'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.70%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.16'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.49%  '
$ws.Range('E22').Value = '  +0.28%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.34'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +6.40%  '
$ws.Range('E24').Value = '  -2.22%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '171.63'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.30%  '
$ws.Range('E26').Value = '  -2.82%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.38'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.05%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.120'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.08%  '
$ws.Range('E29').Value = '  +0.19%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.80'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.18%  '
$ws.Range('E31').Value = '  -0.48%  '
$ws.Range('E32').Value = '  -0.46%  '
$ws.Range('E33').Value = '  -0.31%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.80'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.61%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '88.17'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +8.67%  '
$ws.Range('E36').Value = '  -0.19%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.315.22'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.37%  '
$ws.Range('B39').Value = 'InjectiveProtocol'
$ws.Range('C39').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '14.78'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +12.45%  '
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0187'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.65%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.33'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.50%  '
$ws.Range('E42').Value = '  +5.27%  '
$ws.Range('E43').Value = '  +0.48%  '
$ws.Range('E44').Value = '  +0.07%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.937'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.18%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0519'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +3.93%  '
$ws.Range('E47').Value = '  -0.03%  '
$ws.Range('E48').Value = '  +0.08%  '
$ws.Range('E49').Value = '  +0.26%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '100.72'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.81%  '
$ws.Range('E51').Value = '  +0.59%  '
